$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.015.79'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '1.641.12'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  -0.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.75'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5094'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.76%  '
$ws.Range('E7').Value = '  -0.35%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2562'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06343'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.53'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07764'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.284'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('D13').Value = '1.650.73'
$ws.Range('E13').Value = '  +0.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5429'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '64.08'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.23%  '
$ws.Range('D16').Value = '0.0₅7701'
$ws.Range('E16').Value = '  -2.23%  '
$ws.Range('D17').Value = '26.039.35'
$ws.Range('E17').Value = '  +0.26%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.002'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '198.78'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.420'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.893'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.042'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.98%  '
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.867'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '140.88'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.44%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1192'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.804'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.82%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.58'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.98%  '
$ws.Range('E29').Value = '  -0.94%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.04888'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.71%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.254'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.43%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.159'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.21%  '
$ws.Range('E33').Value = '  -0.65%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.364'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9018'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.82%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.580'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.97%  '
$ws.Range('D37').Value = '1.142.51'
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('E38').Value = '  -2.08%  '
$ws.Range('E39').Value = '  -0.18%  '
$ws.Range('E40').Value = '  -0.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.528'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.54%  '
$ws.Range('D42').Value = '0.0₈128'
$ws.Range('E42').Value = '  +7.80%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8105'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.89%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.34'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.387'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.37%  '
$ws.Range('D46').Value = '1.782.19'
$ws.Range('E46').Value = '  +0.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4528'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.10%  '
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.002'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.65%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.95'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05075'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.001'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.61%  '
